# Add a new "Reference Article" slide at the end of the deck (Title Only layout).
$p = $ppt.ActivePresentation

$s = $p.Slides.Add($p.Slides.Count + 1, 7)  # 7 = ppLayoutTitleOnly

# Title placeholder text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Reference Article"

# Citation rectangle (EMU / 12700 = points, since Shapes.AddShape takes points).
$rect = $s.Shapes.AddShape(1, 80, 163.5295275590551, 583.5294488188977, 72.7031496062992)
$rect.TextFrame.WordWrap = $true
$rect.TextFrame.AutoSize = 1

$tr = $rect.TextFrame.TextRange
$tr.Text = "ieeexplore.ieee.org"
$tr.Font.Bold = $true

$tr = $tr.InsertAfter("/iel5/5888675/5928760/05928901.pdf")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter([char]13 + "Smart Parking Reservation System By ")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter("Mohit")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter(" ")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter("Patil")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter(" and Rahul ")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter("Sakore")
$tr.Font.Bold = $true

$tr = $tr.InsertAfter(" ")
$tr.Font.Bold = $true
